# feat: add 2022-Q3 data
#
# - Insert a new "2022-Q3" sheet right after "总计", holding the latest
#   fund-holding snapshot (copied from "2022-Q2" so it inherits the same
#   look/formatting), then overwrite its values with the new quarter's data.
# - Insert a new row at the top of the "总计" summary table for the new
#   quarter and refresh the sequential index column / trailing row.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2    = $wb.Worksheets.Item("2022-Q2")

# --- 1. Create the new "2022-Q3" sheet right after "总计" -------------------
# Duplicate "2022-Q2" (rather than Worksheets.Add) so the new tab keeps the
# same header row / column formatting, then drop it in place as tab #2 and
# rename it.
$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Force the code/name/ratio columns back to text (matches the source data
# types - e.g. fund code "163302" must stay a string, not become a number).
$wsQ3.Range("B2:G2").NumberFormat = "@"
$wsQ3.Range("B2").Value = "163302"
$wsQ3.Range("C2").Value = "大摩资源优选混合（LOF）"
$wsQ3.Range("D2").Value = "5.08"
$wsQ3.Range("E2").Value = "79.56"
$wsQ3.Range("F2").Value = "3.93"
$wsQ3.Range("G2").Value = "0.1996"
$wsQ3.Range("H2").Value = 4

# --- 2. Update the "总计" roll-up sheet -------------------------------------
# Push the existing rows down one and add the new 2022-Q3 entry at the top.
$wsTotal.Rows("2:2").Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.2

# Keep the sequential index column (A) in order for the rows that were
# pushed down ...
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4

# ... and append the trailing "2020-Q4" row that fell off the bottom.
$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").Value = "2020-Q4"
$wsTotal.Range("C7").Value = 1
$wsTotal.Range("D7").Value = 0.66
